$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Append a new record (Marouane Chemek) as row 46
$row = 46
$ws.Cells.Item($row, 1).Value = "Marouane"
$ws.Cells.Item($row, 2).Value = "Chemek "
$ws.Cells.Item($row, 3).Value = "South Ural State University"
$ws.Cells.Item($row, 4).Value = "Russie"
$ws.Cells.Item($row, 5).Value = "imC8he8AAAAJ"
$ws.Cells.Item($row, 6).Value = "M"
$ws.Cells.Item($row, 7).Value = 1988
$ws.Cells.Item($row, 8).Value = "Médecine, Biologie et Sciences de la Santé"

# Carry over the Genre column's existing look (Arial 8pt) by copying the
# format from the cell directly above, matching how the rest of the column
# was formatted.
$ws.Range("F45").Copy()
$ws.Range("F46").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# Scroll the view down and land the selection on the next empty row, as in
# the saved workbook.
$win = $excel.ActiveWindow
$win.ScrollRow = 27
$ws.Range("H47").Select()
